$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 447.72675691069213
$ws.Range("E2").Value = 182.89167784199395
$ws.Range("F2").Value = 51.188111908216221
$ws.Range("G2").Value = 112.1721279456065

$ws.Range("C3").Value = 32.416364714563123
$ws.Range("D3").Value = 436.99235134094602
$ws.Range("E3").Value = 186.74848824692378
$ws.Range("F3").Value = 48.440894912266103
$ws.Range("G3").Value = 109.39897472307671

$ws.Range("C4").Value = 32.91097173227854
$ws.Range("D4").Value = 441.45440502531841
$ws.Range("E4").Value = 190.07425385812425
$ws.Range("F4").Value = 49.067314862631456
$ws.Range("G4").Value = 110.30597327572663

$ws.Range("C5").Value = 33.308083785147147
$ws.Range("D5").Value = 445.58723892673299
$ws.Range("E5").Value = 192.81465266787015
$ws.Range("F5").Value = 49.629281580298532
$ws.Range("G5").Value = 111.20019906705092
